$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = 9752.14
$ws.Range("B12").Value = 9739.48
$ws.Range("C12").Value = 77.78
$ws.Range("D12").Value = 77.88
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = 0.13

# G12 is a date/time value; copy the format from G11 so the style (numFmt 22) is preserved.
$ws.Range("G11").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G12").Value = 42620.766076388885

$ws.Range("H12").Value = $true
